$d = $word.ActiveDocument

# Locate the word "Mensagem" inside "Enviar Mensagem Lembrete Compromisso"
# (table cell for use case UC007) and replace it with "E-mail", splitting
# the original single run into three runs so the bookmark that sits right
# after the word stays anchored between run 2 ("E-mail") and run 3
# (" Lembrete Compromisso").
$rng = $d.Content.Duplicate
$rng.Find.Execute("Mensagem", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rng.Text = "E-mail"

# Touch formatting on the replaced range and put it straight back so the
# range is materialized as its own run (distinct from its neighbours)
# without altering the final character formatting.
$rng.Bold = 1
$rng.Bold = 0
